$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where "Absent" (column H) becomes 1
$absentRows = @(3,6,7,8,9,10,11,14,15,16,17,18)
foreach ($r in $absentRows) {
    $ws.Range("H$r").Value = 1
}

# Rows where "Total Attendance Count" (D) and "Real" (E) become 1
$presentRows = @(4,5,12,13)
foreach ($r in $presentRows) {
    $ws.Range("D$r").Value = 1
    $ws.Range("E$r").Value = 1
}

# Row 11 also has "Invalid" (G) becoming 1
$ws.Range("G11").Value = 1
